$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4112.272962000001
$ws.Range("D2").Value = 8403.191820900001
$ws.Range("E2").Value = 8563.2335649
$ws.Range("F2").Value = 8678.263568400002
$ws.Range("G2").Value = 8784.291223800001
$ws.Range("H2").Value = 8881.316531100001
$ws.Range("I2").Value = 8971.340012100001
$ws.Range("J2").Value = 9033.356187900001
$ws.Range("K2").Value = 9091.371320100001
$ws.Range("L2").Value = 9146.3856696
$ws.Range("M2").Value = 9197.3989755
$ws.Range("N2").Value = 9245.411498700001
$ws.Range("O2").Value = 9291.423500100002
$ws.Range("P2").Value = 9334.434718800001
$ws.Range("Q2").Value = 9374.445154800002
# Row 3
$ws.Range("C3").Value = 116.2945998
$ws.Range("D3").Value = 237.6239652
$ws.Range("E3").Value = 242.1495972
$ws.Range("F3").Value = 245.4023952
$ws.Range("G3").Value = 248.4006264
$ws.Range("H3").Value = 251.1442908
$ws.Range("I3").Value = 253.6899588
$ws.Range("J3").Value = 255.4436412
$ws.Range("K3").Value = 257.0841828
$ws.Range("L3").Value = 258.6398688
$ws.Range("M3").Value = 260.082414
$ws.Range("N3").Value = 261.4401036
$ws.Range("O3").Value = 262.7412228
$ws.Range("P3").Value = 263.9574864
$ws.Range("Q3").Value = 265.0888944
# Row 4
$ws.Range("C4").Value = 86.9212215
$ws.Range("D4").Value = 178.2179739000001
$ws.Range("E4").Value = 181.6121979
$ws.Range("F4").Value = 184.0517964
$ws.Range("G4").Value = 186.3004698
$ws.Range("H4").Value = 188.3582181
$ws.Range("I4").Value = 190.2674691
$ws.Range("J4").Value = 191.5827309
$ws.Range("K4").Value = 192.8131371
$ws.Range("L4").Value = 193.9799016
$ws.Range("M4").Value = 195.0618105
$ws.Range("N4").Value = 196.0800777
$ws.Range("O4").Value = 197.0559171
$ws.Range("P4").Value = 197.9681148000001
$ws.Range("Q4").Value = 198.8166708000001
# Row 6
$ws.Range("C6").Value = 127432.505286
$ws.Range("D6").Value = 260412.5934090001
$ws.Range("E6").Value = 265372.242849
$ws.Range("F6").Value = 268936.990884
$ws.Range("G6").Value = 272222.758638
$ws.Range("H6").Value = 275229.546111
$ws.Range("I6").Value = 278019.348921
$ws.Range("J6").Value = 279941.2130790001
$ws.Range("K6").Value = 281739.086001
$ws.Range("L6").Value = 283443.965496
$ws.Range("M6").Value = 285024.853755
$ws.Range("N6").Value = 286512.748587
$ws.Range("O6").Value = 287938.647801
$ws.Range("P6").Value = 289271.553588
$ws.Range("Q6").Value = 290511.4659480001
# Row 8
$ws.Range("D8").Value = 7920.000000000001
$ws.Range("E8").Value = 8712.000000000002
$ws.Range("F8").Value = 9583.200000000003
$ws.Range("G8").Value = 10541.52
$ws.Range("H8").Value = 11595.67200000001
$ws.Range("I8").Value = 12755.23920000001
$ws.Range("J8").Value = 14030.76312000001
$ws.Range("K8").Value = 15433.83943200001
$ws.Range("L8").Value = 16977.22337520001
$ws.Range("M8").Value = 18674.94571272002
$ws.Range("N8").Value = 20542.44028399202
$ws.Range("O8").Value = 22596.68431239123
$ws.Range("P8").Value = 24856.35274363035
$ws.Range("Q8").Value = 27341.98801799339
# Row 9
$ws.Range("C9").Value = 2400
$ws.Range("D9").Value = 2760
$ws.Range("E9").Value = 3174
$ws.Range("F9").Value = 3650.099999999999
$ws.Range("G9").Value = 4197.614999999999
$ws.Range("H9").Value = 4827.257249999999
$ws.Range("I9").Value = 5551.345837499998
$ws.Range("J9").Value = 6384.047713124997
$ws.Range("K9").Value = 7341.654870093746
$ws.Range("L9").Value = 8442.903100607808
$ws.Range("M9").Value = 9709.338565698978
$ws.Range("N9").Value = 11165.73935055382
$ws.Range("O9").Value = 12840.6002531369
$ws.Range("P9").Value = 14766.69029110743
$ws.Range("Q9").Value = 16981.69383477354
# Row 12
$ws.Range("C12").Value = 226632.505286
$ws.Range("D12").Value = 271692.5934090001
$ws.Range("E12").Value = 277858.242849
$ws.Range("F12").Value = 282770.290884
$ws.Range("G12").Value = 287561.893638
$ws.Range("H12").Value = 292252.475361
$ws.Range("I12").Value = 296925.9339585
$ws.Range("J12").Value = 300956.0239121251
$ws.Range("K12").Value = 305114.5803030938
$ws.Range("L12").Value = 309464.0919718078
$ws.Range("M12").Value = 314009.138033419
$ws.Range("N12").Value = 318820.9282215459
$ws.Range("O12").Value = 323975.9323665281
$ws.Range("P12").Value = 329494.5966227378
$ws.Range("Q12").Value = 335435.147800767
# Row 14
$ws.Range("C14").Value = 3185.812632150001
$ws.Range("D14").Value = 6510.314835225002
$ws.Range("E14").Value = 6634.306071225
$ws.Range("F14").Value = 6723.424772100001
$ws.Range("G14").Value = 6805.56896595
$ws.Range("H14").Value = 6880.738652775
$ws.Range("I14").Value = 6950.483723025001
$ws.Range("J14").Value = 6998.530326975002
$ws.Range("K14").Value = 7043.477150025001
$ws.Range("L14").Value = 7086.0991374
$ws.Range("M14").Value = 7125.621343875
$ws.Range("N14").Value = 7162.818714675001
$ws.Range("O14").Value = 7198.466195025001
$ws.Range("P14").Value = 7231.788839700001
$ws.Range("Q14").Value = 7262.786648700003
# Row 15
$ws.Range("C15").Value = 46541.818188
$ws.Range("D15").Value = 71340.77743875002
$ws.Range("E15").Value = 72699.48763874998
$ws.Range("F15").Value = 73676.06059500002
$ws.Range("G15").Value = 74576.2061025
$ws.Range("H15").Value = 75399.92416125
$ws.Range("I15").Value = 76164.19864875001
$ws.Range("J15").Value = 76690.69885125001
$ws.Range("K15").Value = 77183.23129875
$ws.Range("L15").Value = 77650.28792999999
$ws.Range("M15").Value = 78083.37680625002
$ws.Range("N15").Value = 78490.98986625001
$ws.Range("O15").Value = 78881.61904875001
$ws.Range("P15").Value = 79246.772415
$ws.Range("Q15").Value = 79586.44996500001
# Row 16
$ws.Range("C16").Value = 1355.637337237114
$ws.Range("D16").Value = 2770.40651743299
$ws.Range("E16").Value = 2823.169884030928
$ws.Range("F16").Value = 2861.093553773197
$ws.Range("G16").Value = 2896.04928414433
$ws.Range("H16").Value = 2928.03707514433
$ws.Range("I16").Value = 2957.71646885567
$ws.Range("J16").Value = 2978.162273412371
$ws.Range("K16").Value = 2997.288993804124
$ws.Range("L16").Value = 3015.426401072165
$ws.Range("M16").Value = 3032.244724175258
$ws.Range("N16").Value = 3048.073734154639
$ws.Range("O16").Value = 3063.243202051547
$ws.Range("P16").Value = 3077.423356824743
$ws.Range("Q16").Value = 3090.614198474228
# Row 19
$ws.Range("C19").Value = 74433.79492638713
$ws.Range("D19").Value = 111973.390729408
$ws.Range("E19").Value = 113508.8555320059
$ws.Range("F19").Value = 114612.4708588732
$ws.Range("G19").Value = 115629.7162905943
$ws.Range("H19").Value = 116560.5918271693
$ws.Range("I19").Value = 117424.2907786307
$ws.Range("J19").Value = 118019.2833896374
$ws.Range("K19").Value = 118575.8893805791
$ws.Range("L19").Value = 119103.7054064722
$ws.Range("M19").Value = 119593.1348123003
$ws.Range("N19").Value = 120053.7742530796
$ws.Range("O19").Value = 120495.2203838266
$ws.Range("P19").Value = 120907.8765495247
$ws.Range("Q19").Value = 121291.7427501742
# Row 20
$ws.Range("C20").Value = 152198.7103596129
$ws.Range("D20").Value = 159719.2026795921
$ws.Range("E20").Value = 164349.3873169941
$ws.Range("F20").Value = 168157.8200251268
$ws.Range("G20").Value = 171932.1773474057
$ws.Range("H20").Value = 175691.8835338307
$ws.Range("I20").Value = 179501.6431798694
$ws.Range("J20").Value = 182936.7405224877
$ws.Range("K20").Value = 186538.6909225147
$ws.Range("L20").Value = 190360.3865653356
$ws.Range("M20").Value = 194416.0032211188
$ws.Range("N20").Value = 198767.1539684662
$ws.Range("O20").Value = 203480.7119827016
$ws.Range("P20").Value = 208586.720073213
$ws.Range("Q20").Value = 214143.4050505928
# Row 22
$ws.Range("D22").Value = 105678
$ws.Range("E22").Value = 105678
$ws.Range("F22").Value = 105678
$ws.Range("G22").Value = 105678
$ws.Range("H22").Value = 105678
$ws.Range("I22").Value = 105678
$ws.Range("J22").Value = 105678
$ws.Range("K22").Value = 105678
$ws.Range("L22").Value = 105678
$ws.Range("M22").Value = 105678
$ws.Range("N22").Value = 105678
$ws.Range("O22").Value = 105678
$ws.Range("P22").Value = 105678
$ws.Range("Q22").Value = 105678
# Row 23
$ws.Range("D23").Value = 8454.24
$ws.Range("E23").Value = 8454.24
$ws.Range("F23").Value = 8454.24
$ws.Range("G23").Value = 8454.24
$ws.Range("H23").Value = 8454.24
$ws.Range("I23").Value = 8454.24
$ws.Range("J23").Value = 8454.24
$ws.Range("K23").Value = 8454.24
$ws.Range("L23").Value = 8454.24
$ws.Range("M23").Value = 8454.24
$ws.Range("N23").Value = 8454.24
$ws.Range("O23").Value = 8454.24
$ws.Range("P23").Value = 8454.24
$ws.Range("Q23").Value = 8454.24
# Row 26
$ws.Range("D26").Value = 140532.24
$ws.Range("E26").Value = 140532.24
$ws.Range("F26").Value = 140532.24
$ws.Range("G26").Value = 140532.24
$ws.Range("H26").Value = 140532.24
$ws.Range("I26").Value = 140532.24
$ws.Range("J26").Value = 140532.24
$ws.Range("K26").Value = 140532.24
$ws.Range("L26").Value = 140532.24
$ws.Range("M26").Value = 140532.24
$ws.Range("N26").Value = 140532.24
$ws.Range("O26").Value = 140532.24
$ws.Range("P26").Value = 140532.24
$ws.Range("Q26").Value = 140532.24
# Row 27
$ws.Range("C27").Value = 76988.71035961289
$ws.Range("D27").Value = 19186.96267959208
$ws.Range("E27").Value = 23817.14731699409
$ws.Range("F27").Value = 27625.58002512681
$ws.Range("G27").Value = 31399.93734740571
$ws.Range("H27").Value = 35159.64353383071
$ws.Range("I27").Value = 38969.40317986938
$ws.Range("J27").Value = 42404.50052248774
$ws.Range("K27").Value = 46006.45092251466
$ws.Range("L27").Value = 49828.14656533563
$ws.Range("M27").Value = 53883.76322111877
$ws.Range("N27").Value = 58234.91396846625
$ws.Range("O27").Value = 62948.47198270159
$ws.Range("P27").Value = 68054.48007321302
$ws.Range("Q27").Value = 73611.1650505928
# Row 32
$ws.Range("C32").Value = 38789.69124601546
$ws.Range("D32").Value = -28720.9476217703
$ws.Range("E32").Value = -24090.7629843683
$ws.Range("F32").Value = -20282.33027623557
$ws.Range("G32").Value = -16507.97295395668
$ws.Range("H32").Value = 15741.86115830081
$ws.Range("I32").Value = 19551.62080433948
$ws.Range("J32").Value = 22986.71814695784
$ws.Range("K32").Value = 26588.66854698476
$ws.Range("L32").Value = 30410.36418980573
$ws.Range("M32").Value = 34465.98084558887
$ws.Range("N32").Value = 38817.13159293635
$ws.Range("O32").Value = 43530.68960717169
$ws.Range("P32").Value = 48636.69769768312
$ws.Range("Q32").Value = 54193.38267506289
# Row 33
$ws.Range("C33").Value = 17.38608901552959
$ws.Range("D33").Value = -12.87313551416203
$ws.Range("E33").Value = -10.79782117990632
$ws.Range("F33").Value = -9.090827699259538
$ws.Range("G33").Value = -7.399107289180288
$ws.Range("H33").Value = 7.055725131517756
$ws.Range("I33").Value = 8.763313364528113
$ws.Range("J33").Value = 10.30297264660365
$ws.Range("K33").Value = 11.91741783223847
$ws.Range("L33").Value = 13.63035594806253
$ws.Range("M33").Value = 15.44814077504419
$ws.Range("N33").Value = 17.39838816767179
$ws.Range("O33").Value = 19.51107163028576
$ws.Range("P33").Value = 21.7996567755662
$ws.Range("Q33").Value = 24.29024168471793
# Row 34
$ws.Range("B34").Value = 6446.184730310095
$ws.Range("C34").Value = -66317.93929336435
$ws.Range("D34").Value = -95038.88691513466
$ws.Range("E34").Value = -119129.649899503
$ws.Range("F34").Value = -139411.9801757385
$ws.Range("G34").Value = -155919.9531296952
$ws.Range("H34").Value = -140178.0919713944
$ws.Range("I34").Value = -120626.4711670549
$ws.Range("J34").Value = -97639.75302009709
$ws.Range("K34").Value = -71051.08447311234
$ws.Range("L34").Value = -40640.72028330661
$ws.Range("M34").Value = -6174.739437717741
$ws.Range("N34").Value = 76173.0817623903
$ws.Range("O34").Value = 124809.7794600734
$ws.Range("P34").Value = 179003.1621351363
$ws.Range("Q34").Value = 179003.1621351363
